# Jan 22 - Input updates
# Remove the obsolete TryDermaFlash/Core and Dr.Denese (x2) test-data rows
# (rows 2-4) from Sheet1; everything below shifts up by three rows and the
# sheet's used range shrinks from A1:E25 to A1:E22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:E4").EntireRow.Delete()

# Reset the active cell back to the top of the sheet (the stale A21
# selection from the deleted rows no longer makes sense once the sheet
# only runs through row 22).
$ws.Cells.Item(1, 1).Select()
